$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Review-point outcomes (column E) - all points resolved, mostly Accepted,
# last one Rejected.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "Accepted"
$ws.Range("E3").Value = "Accepted"
$ws.Range("E4").Value = "Accepted"
$ws.Range("E5").Value = "Accepted"
$ws.Range("E6").Value = "Rejected"

# ---------------------------------------------------------------------------
# Reviewer comments (column G) explaining how each point was addressed.
# Requirement IDs inside the comments are bolded, matching the author's
# formatting.
# ---------------------------------------------------------------------------

$g2 = $ws.Range("G2")
$g2.Value = "Added a new reuirement Req_PO1_DGC_SRS_012_V01 for this point."
$g2.Characters(24, 23).Font.Bold = $true

$g3 = $ws.Range("G3")
$g3.Value = "Changed the statement to be more explicit and clear." + [char]10 + "The requirement cannot be divided into 2 requirements, since it's only concerned with a division by zero operation."

$g4 = $ws.Range("G4")
$g4.Value = "Removed requirement Req_PO1_DGC_SRS_010_V01, and appended the useful parts to requirement Req_PO1_DGC_SRS_007_V01"
$g4.Characters(21, 23).Font.Bold = $true
$g4.Characters(91, 23).Font.Bold = $true

$g5 = $ws.Range("G5")
$g5.Value = "Added the desired previous state of the switch in each requirement."

$g6 = $ws.Range("G6")
$g6.Value = "Renamed the State Machine to Flow Chart and added a new reuirement Req_PO1_DGC_SRS_013_V01 to explicitely state that the software design should follow the described flow chart ."
$g6.Characters(68, 23).Font.Bold = $true

# ---------------------------------------------------------------------------
# Formatting: column D/G now vertically centred; column G content wraps and
# is left aligned, matching the "review point" column; header row follows
# suit (D no longer wraps, G now does).
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).VerticalAlignment = -4108
$ws.Columns.Item(7).VerticalAlignment = -4108

$ws.Range("G1").WrapText = $true
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").Font.Bold = $true

$ws.Range("G2:G6").WrapText = $true
$ws.Range("G2:G6").HorizontalAlignment = -4131

# Drop the stray leftover formatting on the (now unused) G8 cell so it goes
# back to an unstyled, empty cell.
$ws.Range("G8").ClearFormats()

# ---------------------------------------------------------------------------
# Conditional formatting: the "Open"/"Closed" highlighting that used to only
# watch column F now also watches G4.
# ---------------------------------------------------------------------------
$cfOpen1 = $ws.Range("G4").FormatConditions.Add(9, 2, "Open")
$cfOpen1.Text = "Open"
$cfOpen1.Formula1 = 'NOT(ISERROR(SEARCH("Open",G4)))'

$cfClosed = $ws.Range("G4").FormatConditions.Add(9, 2, "Closed")
$cfClosed.Text = "Closed"
$cfClosed.Formula1 = 'NOT(ISERROR(SEARCH("Closed",G4)))'

$cfOpen2 = $ws.Range("G4").FormatConditions.Add(9, 2, "Open")
$cfOpen2.Text = "Open"
$cfOpen2.Formula1 = 'NOT(ISERROR(SEARCH("Open",G4)))'

# ---------------------------------------------------------------------------
# Row heights grew slightly to fit the new wrapped comments.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 57.6
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 57.6

# ---------------------------------------------------------------------------
# Selection moved to D13 (what was on screen when the file was last saved).
# ---------------------------------------------------------------------------
$ws.Range("D13").Select()
